$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell as TEXT (not Number),
# without leaving any NumberFormat / quote-prefix style behind. We build the
# text in an off-grid helper cell via a text formula (="...") which yields a
# String-typed result, copy it, and paste-special only the Values into the
# destination - the destination keeps its original (default) cell style and
# ends up with a plain text value, exactly like the source data.
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range("ZZ1").Formula = '="' + $escaped + '"'
    $ws.Range("ZZ1").Copy()
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $ws.Range("ZZ1").ClearContents()
}

$ws.Range("D2").Value = "27.189.92"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.850.36"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  -0.34%  "
Set-TextValue $ws.Range("D5") "313.58"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  -0.26%  "
Set-TextValue $ws.Range("D7") "0.4601"
$ws.Range("E7").Value = "  -0.18%  "
Set-TextValue $ws.Range("D8") "0.3700"
$ws.Range("E8").Value = "  -0.22%  "
Set-TextValue $ws.Range("D9") "0.07273"
$ws.Range("E9").Value = "  -0.96%  "
Set-TextValue $ws.Range("D10") "0.8827"
$ws.Range("E10").Value = "  +0.80%  "
Set-TextValue $ws.Range("D11") "19.99"
$ws.Range("E11").Value = "  +0.75%  "
Set-TextValue $ws.Range("D12") "0.07822"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "1.829.74"
$ws.Range("E13").Value = "  +5.84%  "
Set-TextValue $ws.Range("D14") "5.366"
$ws.Range("E14").Value = "  +0.39%  "
Set-TextValue $ws.Range("D15") "6.492"
$ws.Range("E15").Value = "  -1.20%  "
Set-TextValue $ws.Range("D16") "91.26"
$ws.Range("E16").Value = "  -0.33%  "
Set-TextValue $ws.Range("D17") "1.002"
$ws.Range("E17").Value = "  -0.38%  "
Set-TextValue $ws.Range("D18") "0.000008908"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  -0.16%  "
Set-TextValue $ws.Range("D20") "14.69"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "27.225.24"
$ws.Range("E21").Value = "  +0.45%  "
Set-TextValue $ws.Range("D22") "5.088"
$ws.Range("E22").Value = "  -0.63%  "
Set-TextValue $ws.Range("D23") "10.50"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").Value = "2.072.60"
$ws.Range("E24").Value = "  +0.30%  "
Set-TextValue $ws.Range("D25") "1.931"
$ws.Range("E25").Value = "  +4.87%  "
Set-TextValue $ws.Range("D26") "151.53"
$ws.Range("E26").Value = "  -1.10%  "
Set-TextValue $ws.Range("D27") "18.34"
$ws.Range("E27").Value = "  -0.46%  "
Set-TextValue $ws.Range("D28") "2.056"
$ws.Range("E28").Value = "  +0.31%  "
Set-TextValue $ws.Range("D29") "115.53"
$ws.Range("E29").Value = "  -0.08%  "
Set-TextValue $ws.Range("D30") "5.039"
$ws.Range("E30").Value = "  -2.37%  "
Set-TextValue $ws.Range("D31") "0.08815"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  +4.43%  "
Set-TextValue $ws.Range("D33") "0.7591"
$ws.Range("E33").Value = "  +3.48%  "
Set-TextValue $ws.Range("D34") "1.167"
$ws.Range("E34").Value = "  +3.17%  "
Set-TextValue $ws.Range("D35") "4.493"
$ws.Range("E35").Value = "  +1.36%  "
Set-TextValue $ws.Range("D36") "2.711"
$ws.Range("E36").Value = "  +8.88%  "
Set-TextValue $ws.Range("D37") "1.083"
$ws.Range("E37").Value = "  +1.09%  "
Set-TextValue $ws.Range("D38") "0.01942"
$ws.Range("E38").Value = "  -0.51%  "
Set-TextValue $ws.Range("D39") "0.05229"
$ws.Range("E39").Value = "  -0.19%  "
Set-TextValue $ws.Range("D40") "2.945"
$ws.Range("E40").Value = "  +0.08%  "
Set-TextValue $ws.Range("D41") "7.047"
$ws.Range("E41").Value = "  -1.14%  "
Set-TextValue $ws.Range("D42") "0.5093"
$ws.Range("E42").Value = "  -1.39%  "
Set-TextValue $ws.Range("D43") "0.1622"
$ws.Range("E43").Value = "  -0.35%  "
Set-TextValue $ws.Range("D44") "8.358"
$ws.Range("E44").Value = "  +1.67%  "
Set-TextValue $ws.Range("D45") "0.4773"
$ws.Range("E45").Value = "  -1.42%  "
Set-TextValue $ws.Range("D46") "10.33"
$ws.Range("E46").Value = "  +0.98%  "
Set-TextValue $ws.Range("D47") "1.002"
$ws.Range("E47").Value = "  -0.33%  "
Set-TextValue $ws.Range("D48") "102.36"
$ws.Range("E48").Value = "  +0.06%  "
Set-TextValue $ws.Range("D49") "1.633"
$ws.Range("E49").Value = "  +0.07%  "
Set-TextValue $ws.Range("D50") "0.06210"
$ws.Range("E50").Value = "  +0.15%  "
Set-TextValue $ws.Range("D51") "65.55"
$ws.Range("E51").Value = "  +1.62%  "
